$d = $word.ActiveDocument

$pairs = @(
    @("11×40=", "71×25="),
    @("26×82=", "75×38="),
    @("16×65=", "40×39="),
    @("29×70=", "77×63="),
    @("84×35=", "89×37="),
    @("66×77=", "36×22="),
    @("11×97=", "74×71="),
    @("43×91=", "52×91="),
    @("52×77=", "48×22="),
    @("81×21=", "59×73="),
    @("63×16=", "44×82="),
    @("83×44=", "57×85="),
    @("67×71=", "69×39="),
    @("44×20=", "26×68="),
    @("40×91=", "12×78="),
    @("52×82=", "29×89="),
    @("51×15=", "32×30="),
    @("86×19=", "90×37="),
    @("61×94=", "36×84="),
    @("85×72=", "63×94="),
    @("38×61=", "19×18="),
    @("80×17=", "38×69="),
    @("65×13=", "71×98="),
    @("89×25=", "52×76="),
    @("96×69=", "36×85=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
